$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The registration that used to occupy row 2 (Turai Attila / 980) is
# gone, so every later entry moves up one row.
$ws.Rows.Item(2).Delete()

# The stray placeholder row that used to be row 8 (bare "982" id, no
# other data) shifted up to row 7 - drop it too, leaving the sheet
# ending at row 6.
$ws.Rows.Item(7).Delete()

# What is now row 6 (formerly "Gál László") is not a real registration
# for this competition - blank out the registrant columns but keep the
# row.
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Automatikus verseny ID töltés: stamp every remaining registrant with
# the competition id.
$ws.Range("V2:V6").Value = "VID_00001"
